$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix tiny precision glitch on row 4 timestamps
$ws.Range("B4").Value = 44327.18941920139
$ws.Range("C4").Value = 44327.2041735301

# Row 5
$ws.Range("A5").Value = "Training"
$ws.Range("B5").Value = 44341.75635722222
$ws.Range("C5").Value = 44341.77731420139
$ws.Range("B5").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("C5").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("D5").Value = 3
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 3
$ws.Range("G5").Value = 28
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 26
$ws.Range("J5").Value = 2
$ws.Range("K5").Value = 1

# Row 6
$ws.Range("A6").Value = "Training"
$ws.Range("B6").Value = 44341.79840987932
$ws.Range("C6").Value = 44341.81956005171
$ws.Range("B6").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("C6").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("D6").Value = 31
$ws.Range("E6").Value = 7
$ws.Range("F6").Value = 3
$ws.Range("G6").Value = 4
$ws.Range("H6").Value = 3
$ws.Range("I6").Value = 25
$ws.Range("J6").Value = 3
$ws.Range("K6").Value = 0
